$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $old"
    }
}

# 1. "it is important that verify the site's credibility" -> "... that you verify ..."
Replace-Text "it is important that verify the site" "it is important that you verify the site"

# 2. "phone number to trick people" -> "phone numbers to trick people"
Replace-Text "s and phone number to trick people" "s and phone numbers to trick people"

# 3. Remove the empty paragraph that sits between the "...could be a scam." paragraph
#    and the "When conduct my personal audit..." paragraph.
for ($i = 2; $i -lt $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Length -le 1) {
        $prev = $d.Paragraphs.Item($i - 1)
        $next = $d.Paragraphs.Item($i + 1)
        if ($prev.Range.Text -like "*it could be a scam.*" -and $next.Range.Text -like "When conduct my personal audit*") {
            $p.Range.Delete()
            break
        }
    }
}

# 4. "When conduct my personal audit, ... asking that if I wanted " ->
#    "When I conducted my personal audit, ... school email, I was constantly getting ... wanted a "
Replace-Text "When conduct my personal audit, I noticed that with my school email I was constant getting emails from unknown sources asking that if I wanted " "When I conducted my personal audit, I noticed that with my school email, I was constantly getting emails from unknown sources asking that if I wanted a "

# 5. "this was clea" -> "it was clea"
Replace-Text "this was clea" "it was clea"

# 6. Passwords paragraph: several small wording fixes.
Replace-Text "need to be changed" "needed to be changed"
Replace-Text "strong were not" "strong, were not"
Replace-Text "checked them use a" "checked them using a"

# 7. "I enable two factor" -> "I enabled two factor"
Replace-Text "I enable two factor" "I enabled two factor"

# 8. "password, the step" -> "password, the extra step"; "Duo mobile" -> "Duo Mobile"; "logging into." -> "logging in."
Replace-Text "password, the s" "password, the extra s"
Replace-Text "Duo mobile" "Duo Mobile"
Replace-Text "logging into." "logging in."

# 9. "logged on a new device" -> "logged into a new device"
Replace-Text "logged on a new device" "logged into a new device"
